$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column A; everything currently in A:AD shifts to B:AE.
$ws.Range("A1").EntireColumn.Insert()

# New header cell for the inserted "In Stash:" column.
$ws.Range("A1").Value = "In Stash:"

# Populate the new column with the per-part stash marker. Assign values in an
# order that reproduces the author's original sharedStrings insertion order
# (In Stash:, +, THONK, *) rather than the sheet's row order.
$ws.Range("A6").Value = "+"
$ws.Range("A7").Value = "THONK"
$ws.Range("A4").Value = "*"

$ws.Range("A5").Value = "THONK"
$ws.Range("A8").Value = "+"
$ws.Range("A9").Value = "+"
$ws.Range("A10").Value = "*"

$ws.Range("A15").Value = "+"
$ws.Range("A16").Value = "*"
$ws.Range("A17").Value = "*"
$ws.Range("A18").Value = "*"
$ws.Range("A19").Value = "*"
$ws.Range("A20").Value = "*"
$ws.Range("A21").Value = "+"
$ws.Range("A22").Value = "*"
$ws.Range("A23").Value = "+"
$ws.Range("A24").Value = "+"
$ws.Range("A25").Value = "+"
$ws.Range("A26").Value = "+"
$ws.Range("A27").Value = "+"

# Fix the sort range / condition (EntireColumn.Insert does not reflow sortState).
$ws.Range("B15:AE27").Sort($ws.Range("F15"))

$ws.Range("E12").Select()
